$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PA3 and PA4 deadlines (pushed back to add pa4 instructions and links)
$ws.Range("B4").Value = "September 27, 2024"
$ws.Range("B5").Value = "October 04, 2024"

# Update the active selection to B5
$ws.Range("B5").Select()
